# Release pcmt-vax-prequal 0.2.0
# Updates the Metadata sheet (Version, Status, Date, FHIR Version) and
# removes the "Mapping: RIM Mapping" column from the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.2.0"
$meta.Range("B6").Value = "active"
$meta.Range("B8").Value = "2025-09-16T20:42:07+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet: drop the "Mapping: RIM Mapping" column (AK) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AK1").EntireColumn.Delete()
